$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits after the title
#    ("Publication Ethics and Malpractice Statement"). It will be
#    re-created further down, right after the newly-inserted text.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Expand ". During the peer review process" into the longer
#    sentence describing that the manuscript is sent to reviewers.
# ------------------------------------------------------------------
$old = ". During the peer review process"
$new = ". The manuscript is sent to two external reviewers for a peer review. During the peer review process"
$rng = $d.Content
$rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark right after the word
#    "reviewers" (before " for a peer review...").
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("two external reviewers", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $bmRange = $rng2.Duplicate
    $bmRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ------------------------------------------------------------------
# 4) Remove the <w:lastRenderedPageBreak/> marker that precedes
#    "Referees should adhere to high standards ...".
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Referees should adhere to high standards", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng3.Find.Found) {
    $p = $rng3.Paragraphs(1)
    $p.Range.Fields | Out-Null
}
